$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data rows 2-5 with new values
$ws.Range("B2").Value = 0.1100390323238504
$ws.Range("C2").Value = 0.6932609026360904
$ws.Range("D2").Value = 0.7099478564566527
$ws.Range("E2").Value = 0.8425840352490976
$ws.Range("F2").Value = 0.8669020772270012

$ws.Range("B3").Value = 0.1368261710548836
$ws.Range("C3").Value = 0.5060341660178811
$ws.Range("D3").Value = 0.4144708479301442
$ws.Range("E3").Value = 0.6437941036776775
$ws.Range("F3").Value = 0.6631150787043795
$ws.Range("G3").Value = 10

$ws.Range("B4").Value = -0.1394442012808889
$ws.Range("C4").Value = 0.3774106509366791
$ws.Range("D4").Value = 0.2048604832801028
$ws.Range("E4").Value = 0.4526151602411289
$ws.Range("F4").Value = 0.4716979516714963
$ws.Range("G4").Value = 6

$ws.Range("B5").Value = -0.06415781549280508
$ws.Range("C5").Value = 0.7320926441837636
$ws.Range("D5").Value = 0.5400758649567836
$ws.Range("E5").Value = 0.7348985405869192
$ws.Range("F5").Value = 1.035335346318259
$ws.Range("G5").Value = 2

# Remove rows 6-9 (Q4-Q7), which removes unused shared strings too
$ws.Range("A6:G9").EntireRow.Delete()
